$d = $word.ActiveDocument

$d.Content.Find.Execute("50×45=2250", $true, $false, $false, $false, $false, $true, 1, $false, "56×51=2856", 2) | Out-Null
$d.Content.Find.Execute("58×76=4408", $true, $false, $false, $false, $false, $true, 1, $false, "89×52=4628", 2) | Out-Null
$d.Content.Find.Execute("18×100=1800", $true, $false, $false, $false, $false, $true, 1, $false, "60×38=2280", 2) | Out-Null
$d.Content.Find.Execute("16×74=1184", $true, $false, $false, $false, $false, $true, 1, $false, "29×72=2088", 2) | Out-Null
$d.Content.Find.Execute("43×52=2236", $true, $false, $false, $false, $false, $true, 1, $false, "68×96=6528", 2) | Out-Null
$d.Content.Find.Execute("77×70=5390", $true, $false, $false, $false, $false, $true, 1, $false, "61×79=4819", 2) | Out-Null
$d.Content.Find.Execute("93×65=6045", $true, $false, $false, $false, $false, $true, 1, $false, "31×70=2170", 2) | Out-Null
$d.Content.Find.Execute("80×10=800", $true, $false, $false, $false, $false, $true, 1, $false, "53×11=583", 2) | Out-Null
$d.Content.Find.Execute("95×61=5795", $true, $false, $false, $false, $false, $true, 1, $false, "83×21=1743", 2) | Out-Null
$d.Content.Find.Execute("83×76=6308", $true, $false, $false, $false, $false, $true, 1, $false, "43×75=3225", 2) | Out-Null
$d.Content.Find.Execute("40×41=1640", $true, $false, $false, $false, $false, $true, 1, $false, "11×94=1034", 2) | Out-Null
$d.Content.Find.Execute("51×56=2856", $true, $false, $false, $false, $false, $true, 1, $false, "21×39=819", 2) | Out-Null
$d.Content.Find.Execute("10×78=780", $true, $false, $false, $false, $false, $true, 1, $false, "13×61=793", 2) | Out-Null
$d.Content.Find.Execute("99×18=1782", $true, $false, $false, $false, $false, $true, 1, $false, "42×11=462", 2) | Out-Null
$d.Content.Find.Execute("12×82=984", $true, $false, $false, $false, $false, $true, 1, $false, "20×24=480", 2) | Out-Null
$d.Content.Find.Execute("92×10=920", $true, $false, $false, $false, $false, $true, 1, $false, "66×68=4488", 2) | Out-Null
$d.Content.Find.Execute("55×44=2420", $true, $false, $false, $false, $false, $true, 1, $false, "41×44=1804", 2) | Out-Null
$d.Content.Find.Execute("97×32=3104", $true, $false, $false, $false, $false, $true, 1, $false, "36×35=1260", 2) | Out-Null
$d.Content.Find.Execute("90×98=8820", $true, $false, $false, $false, $false, $true, 1, $false, "91×38=3458", 2) | Out-Null
$d.Content.Find.Execute("47×25=1175", $true, $false, $false, $false, $false, $true, 1, $false, "58×37=2146", 2) | Out-Null
$d.Content.Find.Execute("36×71=2556", $true, $false, $false, $false, $false, $true, 1, $false, "91×34=3094", 2) | Out-Null
$d.Content.Find.Execute("36×23=828", $true, $false, $false, $false, $false, $true, 1, $false, "30×89=2670", 2) | Out-Null
$d.Content.Find.Execute("75×40=3000", $true, $false, $false, $false, $false, $true, 1, $false, "84×68=5712", 2) | Out-Null
$d.Content.Find.Execute("11×18=198", $true, $false, $false, $false, $false, $true, 1, $false, "18×78=1404", 2) | Out-Null
$d.Content.Find.Execute("26×59=1534", $true, $false, $false, $false, $false, $true, 1, $false, "32×32=1024", 2) | Out-Null
$d.Content.Find.Execute("74×31=2294", $true, $false, $false, $false, $false, $true, 1, $false, "33×80=2640", 2) | Out-Null
$d.Content.Find.Execute("63×15=945", $true, $false, $false, $false, $false, $true, 1, $false, "98×39=3822", 2) | Out-Null
$d.Content.Find.Execute("94×98=9212", $true, $false, $false, $false, $false, $true, 1, $false, "39×96=3744", 2) | Out-Null
$d.Content.Find.Execute("85×15=1275", $true, $false, $false, $false, $false, $true, 1, $false, "36×57=2052", 2) | Out-Null
$d.Content.Find.Execute("92×77=7084", $true, $false, $false, $false, $false, $true, 1, $false, "24×57=1368", 2) | Out-Null
$d.Content.Find.Execute("16×47=752", $true, $false, $false, $false, $false, $true, 1, $false, "11×58=638", 2) | Out-Null
$d.Content.Find.Execute("23×13=299", $true, $false, $false, $false, $false, $true, 1, $false, "47×96=4512", 2) | Out-Null
$d.Content.Find.Execute("49×22=1078", $true, $false, $false, $false, $false, $true, 1, $false, "68×11=748", 2) | Out-Null
$d.Content.Find.Execute("55×40=2200", $true, $false, $false, $false, $false, $true, 1, $false, "64×79=5056", 2) | Out-Null
$d.Content.Find.Execute("81×60=4860", $true, $false, $false, $false, $false, $true, 1, $false, "78×60=4680", 2) | Out-Null
$d.Content.Find.Execute("30×81=2430", $true, $false, $false, $false, $false, $true, 1, $false, "46×36=1656", 2) | Out-Null
$d.Content.Find.Execute("33×15=495", $true, $false, $false, $false, $false, $true, 1, $false, "31×36=1116", 2) | Out-Null
$d.Content.Find.Execute("72×73=5256", $true, $false, $false, $false, $false, $true, 1, $false, "78×26=2028", 2) | Out-Null
$d.Content.Find.Execute("94×60=5640", $true, $false, $false, $false, $false, $true, 1, $false, "52×89=4628", 2) | Out-Null
$d.Content.Find.Execute("11×39=429", $true, $false, $false, $false, $false, $true, 1, $false, "54×77=4158", 2) | Out-Null
$d.Content.Find.Execute("52×61=3172", $true, $false, $false, $false, $false, $true, 1, $false, "73×29=2117", 2) | Out-Null
$d.Content.Find.Execute("89×78=6942", $true, $false, $false, $false, $false, $true, 1, $false, "93×87=8091", 2) | Out-Null
$d.Content.Find.Execute("62×80=4960", $true, $false, $false, $false, $false, $true, 1, $false, "69×32=2208", 2) | Out-Null
$d.Content.Find.Execute("49×64=3136", $true, $false, $false, $false, $false, $true, 1, $false, "14×37=518", 2) | Out-Null
$d.Content.Find.Execute("26×15=390", $true, $false, $false, $false, $false, $true, 1, $false, "100×87=8700", 2) | Out-Null
$d.Content.Find.Execute("60×37=2220", $true, $false, $false, $false, $false, $true, 1, $false, "80×86=6880", 2) | Out-Null
$d.Content.Find.Execute("38×72=2736", $true, $false, $false, $false, $false, $true, 1, $false, "37×49=1813", 2) | Out-Null
$d.Content.Find.Execute("85×28=2380", $true, $false, $false, $false, $false, $true, 1, $false, "99×44=4356", 2) | Out-Null
$d.Content.Find.Execute("90×43=3870", $true, $false, $false, $false, $false, $true, 1, $false, "67×69=4623", 2) | Out-Null
$d.Content.Find.Execute("94×44=4136", $true, $false, $false, $false, $false, $true, 1, $false, "89×22=1958", 2) | Out-Null
$d.Content.Find.Execute("41×26=1066", $true, $false, $false, $false, $false, $true, 1, $false, "58×69=4002", 2) | Out-Null
$d.Content.Find.Execute("65×43=2795", $true, $false, $false, $false, $false, $true, 1, $false, "76×50=3800", 2) | Out-Null
$d.Content.Find.Execute("54×56=3024", $true, $false, $false, $false, $false, $true, 1, $false, "41×78=3198", 2) | Out-Null
$d.Content.Find.Execute("58×89=5162", $true, $false, $false, $false, $false, $true, 1, $false, "66×84=5544", 2) | Out-Null
$d.Content.Find.Execute("17×13=221", $true, $false, $false, $false, $false, $true, 1, $false, "74×82=6068", 2) | Out-Null
$d.Content.Find.Execute("19×61=1159", $true, $false, $false, $false, $false, $true, 1, $false, "48×44=2112", 2) | Out-Null
$d.Content.Find.Execute("28×58=1624", $true, $false, $false, $false, $false, $true, 1, $false, "98×39=3822", 2) | Out-Null
$d.Content.Find.Execute("12×60=720", $true, $false, $false, $false, $false, $true, 1, $false, "25×80=2000", 2) | Out-Null
$d.Content.Find.Execute("83×15=1245", $true, $false, $false, $false, $false, $true, 1, $false, "94×80=7520", 2) | Out-Null
$d.Content.Find.Execute("15×78=1170", $true, $false, $false, $false, $false, $true, 1, $false, "100×69=6900", 2) | Out-Null
$d.Content.Find.Execute("74×21=1554", $true, $false, $false, $false, $false, $true, 1, $false, "77×52=4004", 2) | Out-Null
$d.Content.Find.Execute("56×95=5320", $true, $false, $false, $false, $false, $true, 1, $false, "21×51=1071", 2) | Out-Null
$d.Content.Find.Execute("93×37=3441", $true, $false, $false, $false, $false, $true, 1, $false, "76×98=7448", 2) | Out-Null
$d.Content.Find.Execute("94×89=8366", $true, $false, $false, $false, $false, $true, 1, $false, "53×56=2968", 2) | Out-Null
$d.Content.Find.Execute("56×33=1848", $true, $false, $false, $false, $false, $true, 1, $false, "49×18=882", 2) | Out-Null
$d.Content.Find.Execute("63×72=4536", $true, $false, $false, $false, $false, $true, 1, $false, "41×19=779", 2) | Out-Null
$d.Content.Find.Execute("67×20=1340", $true, $false, $false, $false, $false, $true, 1, $false, "37×61=2257", 2) | Out-Null
$d.Content.Find.Execute("76×51=3876", $true, $false, $false, $false, $false, $true, 1, $false, "48×46=2208", 2) | Out-Null
$d.Content.Find.Execute("43×20=860", $true, $false, $false, $false, $false, $true, 1, $false, "57×64=3648", 2) | Out-Null
$d.Content.Find.Execute("58×42=2436", $true, $false, $false, $false, $false, $true, 1, $false, "56×77=4312", 2) | Out-Null
$d.Content.Find.Execute("38×33=1254", $true, $false, $false, $false, $false, $true, 1, $false, "11×32=352", 2) | Out-Null
$d.Content.Find.Execute("35×85=2975", $true, $false, $false, $false, $false, $true, 1, $false, "52×55=2860", 2) | Out-Null
$d.Content.Find.Execute("58×10=580", $true, $false, $false, $false, $false, $true, 1, $false, "82×37=3034", 2) | Out-Null
$d.Content.Find.Execute("67×29=1943", $true, $false, $false, $false, $false, $true, 1, $false, "60×53=3180", 2) | Out-Null
$d.Content.Find.Execute("32×44=1408", $true, $false, $false, $false, $false, $true, 1, $false, "94×64=6016", 2) | Out-Null
$d.Content.Find.Execute("44×29=1276", $true, $false, $false, $false, $false, $true, 1, $false, "30×63=1890", 2) | Out-Null
$d.Content.Find.Execute("67×63=4221", $true, $false, $false, $false, $false, $true, 1, $false, "12×18=216", 2) | Out-Null
$d.Content.Find.Execute("72×85=6120", $true, $false, $false, $false, $false, $true, 1, $false, "96×45=4320", 2) | Out-Null
$d.Content.Find.Execute("85×22=1870", $true, $false, $false, $false, $false, $true, 1, $false, "67×44=2948", 2) | Out-Null
$d.Content.Find.Execute("93×68=6324", $true, $false, $false, $false, $false, $true, 1, $false, "17×68=1156", 2) | Out-Null
$d.Content.Find.Execute("77×19=1463", $true, $false, $false, $false, $false, $true, 1, $false, "37×26=962", 2) | Out-Null
$d.Content.Find.Execute("83×32=2656", $true, $false, $false, $false, $false, $true, 1, $false, "83×66=5478", 2) | Out-Null
$d.Content.Find.Execute("29×84=2436", $true, $false, $false, $false, $false, $true, 1, $false, "12×55=660", 2) | Out-Null
$d.Content.Find.Execute("88×56=4928", $true, $false, $false, $false, $false, $true, 1, $false, "64×80=5120", 2) | Out-Null
$d.Content.Find.Execute("94×53=4982", $true, $false, $false, $false, $false, $true, 1, $false, "94×10=940", 2) | Out-Null
$d.Content.Find.Execute("92×17=1564", $true, $false, $false, $false, $false, $true, 1, $false, "99×33=3267", 2) | Out-Null
$d.Content.Find.Execute("76×70=5320", $true, $false, $false, $false, $false, $true, 1, $false, "69×100=6900", 2) | Out-Null
$d.Content.Find.Execute("46×90=4140", $true, $false, $false, $false, $false, $true, 1, $false, "76×65=4940", 2) | Out-Null
$d.Content.Find.Execute("32×90=2880", $true, $false, $false, $false, $false, $true, 1, $false, "85×85=7225", 2) | Out-Null
$d.Content.Find.Execute("100×56=5600", $true, $false, $false, $false, $false, $true, 1, $false, "58×93=5394", 2) | Out-Null
$d.Content.Find.Execute("96×53=5088", $true, $false, $false, $false, $false, $true, 1, $false, "46×38=1748", 2) | Out-Null
$d.Content.Find.Execute("38×53=2014", $true, $false, $false, $false, $false, $true, 1, $false, "39×49=1911", 2) | Out-Null
$d.Content.Find.Execute("76×83=6308", $true, $false, $false, $false, $false, $true, 1, $false, "93×51=4743", 2) | Out-Null
$d.Content.Find.Execute("31×76=2356", $true, $false, $false, $false, $false, $true, 1, $false, "35×46=1610", 2) | Out-Null
$d.Content.Find.Execute("54×74=3996", $true, $false, $false, $false, $false, $true, 1, $false, "60×20=1200", 2) | Out-Null
$d.Content.Find.Execute("17×94=1598", $true, $false, $false, $false, $false, $true, 1, $false, "11×52=572", 2) | Out-Null
$d.Content.Find.Execute("77×14=1078", $true, $false, $false, $false, $false, $true, 1, $false, "39×33=1287", 2) | Out-Null
$d.Content.Find.Execute("66×49=3234", $true, $false, $false, $false, $false, $true, 1, $false, "100×14=1400", 2) | Out-Null
$d.Content.Find.Execute("57×17=969", $true, $false, $false, $false, $false, $true, 1, $false, "17×95=1615", 2) | Out-Null
$d.Content.Find.Execute("52×39=2028", $true, $false, $false, $false, $false, $true, 1, $false, "52×78=4056", 2) | Out-Null

Write-Output "Replaced 100 multiplication answers"
